# Add 2022-Q4 data:
#  1. Insert a new "2022-Q4" row at the top of the "总计" (total) summary sheet,
#     pushing the existing quarterly summary rows down by one.
#  2. Insert a brand-new "2022-Q4" worksheet (with the per-fund holdings detail)
#     right after the "总计" sheet, pushing the other quarterly sheets down.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) "总计" sheet: shift existing data rows (2-5) down to (3-6), then write
#    the new 2022-Q4 row into row 2. Use Value2 for reads (Value getter is
#    unreliable in this host) and Value for writes.
# ---------------------------------------------------------------------------
for ($r = 5; $r -ge 2; $r--) {
    $totalSheet.Cells.Item($r + 1, 1).Value = $totalSheet.Cells.Item($r, 1).Value2
    $totalSheet.Cells.Item($r + 1, 2).Value = $totalSheet.Cells.Item($r, 2).Value2
    $totalSheet.Cells.Item($r + 1, 3).Value = $totalSheet.Cells.Item($r, 3).Value2
    $totalSheet.Cells.Item($r + 1, 4).Value = $totalSheet.Cells.Item($r, 4).Value2
}

# Row 6 is brand new - copy column A's number/border/bold style (s="2") from
# the row above it so it matches the other index cells in the column.
$totalSheet.Cells.Item(5, 1).Copy()
$totalSheet.Cells.Item(6, 1).PasteSpecial(-4122)

# New 2022-Q4 summary row
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 0.12

# Refresh the sequential index column (A) for every data row
for ($r = 2; $r -le 6; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2) New "2022-Q4" worksheet with the fund-holdings detail table, inserted
#    right after "总计" (so it becomes the second tab).
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"
$q4Sheet.Outline.SummaryRow = 1
$q4Sheet.Outline.SummaryColumn = 1

# Header row (bold + bordered style, like the other sheets) - copy the style
# from the 总计 sheet's header and then set this sheet's own header text.
$totalSheet.Cells.Item(1, 2).Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 2; $c -le 8; $c++) {
    $q4Sheet.Cells.Item(1, $c).Value = $headers[$c - 2]
}

# Column A index cells also use the bold/bordered style (s="2")
$totalSheet.Cells.Item(2, 1).Copy()
$q4Sheet.Range("A2:A7").PasteSpecial(-4122)

$rows = @(
    @("002307", "银华多元视野灵活配置混合", "1.52", "89.13", "2.12", "0.0322", 7),
    @("002863", "金信深圳成长灵活配置混合", "0.73", "91.56", "3.58", "0.0261", 9),
    @("180028", "银华永祥灵活配置混合", "0.70", "77.51", "3.55", "0.0248", 10),
    @("005117", "金信价值精选灵活配置混合A", "0.76", "92.96", "2.79", "0.0212", 6),
    @("005251", "银华多元动力灵活配置混合", "0.43", "88.30", "3.20", "0.0138", 3),
    @("005118", "金信价值精选灵活配置混合C", "0.06", "92.96", "2.79", "0.0017", 6)
)

$r = 2
foreach ($row in $rows) {
    $q4Sheet.Cells.Item($r, 1).Value = $r - 2
    $q4Sheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $q4Sheet.Cells.Item($r, 3).Value = $row[1]
    $q4Sheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $q4Sheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $q4Sheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $q4Sheet.Cells.Item($r, 7).Value = "'" + $row[5]
    $q4Sheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}
